$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the phone number values in column A
$ws.Range("A3").Value = 7897798797989
$ws.Range("A5").Value = 1234567890

# Update the selected cell/range on the sheet
$ws.Range("B6").Select()
